# Artefato 19 - Regras de Negocios: adiciona RN-008 e RN-009 (eventos 12 e 13 do
# artefato 17, agora classificados como "externo-previsivel") logo apos a RN-007
# existente, e realoca o bookmark oculto _GoBack para depois da nova RN-009.
#
# A RN-007 original guardava o texto em duas runs ("RN-007" e ": O provedor...")
# com o bookmark _GoBack encaixado entre elas. No alvo essas runs viram uma unica
# run com o texto completo, o rPr do run e copiado para o pPr do paragrafo, e o
# bookmark _GoBack passa a viver sozinho num paragrafo novo, logo apos a RN-009.

$d = $word.ActiveDocument

# Localiza o paragrafo "RN-007" dinamicamente em vez de fixar o indice.
$rn007Paragraph = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "RN-007*") {
        $rn007Paragraph = $candidate
        break
    }
}

# Passo 1: reescreve o paragrafo da RN-007 com uma unica run (texto unificado),
# o rPr promovido para o pPr, e sem o bookmark _GoBack (que sera reinserido depois).
$rn007Xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="160" w:afterAutospacing="0"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/></w:rPr><w:t>RN-007: O provedor de cartão de crédito tem até seis minutos para enviar a resposta, se não a ação é cancelada.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rn007Paragraph.Range.InsertXML($rn007Xml)

# Passo 2: insere RN-008, RN-009 (quatro runs: "R" / "N-009" / ":" / " O cliente...")
# e o paragrafo vazio com o bookmark _GoBack, logo antes do paragrafo final vazio
# que ja existia no documento.
$trailingEmptyParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertionPoint = $trailingEmptyParagraph.Range
$insertionPoint.Collapse(1)

$newParagraphsXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="160" w:afterAutospacing="0"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/></w:rPr><w:t>RN-008: O cliente tem o prazo de até 24 horas para efetuar o pagamento em dinheiro, após a conclusão do serviço.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="160" w:afterAutospacing="0"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/></w:rPr><w:t>R</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/></w:rPr><w:t>N-009</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> O cliente após a conclusão de serviço deve efetuar o pagamento.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="160" w:afterAutospacing="0"/></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($newParagraphsXml)

Write-Host "RN-008 e RN-009 inseridas; bookmark _GoBack realocado."
